# Read data from Excel config flexibly for desirable years
#
# Adds a "year" column (with sample years 2020, 2030, 2040) to the
# "config" sheet so downstream code can read the desired years directly
# from the workbook, instead of relying on hard-coded defaults.

$wb = $excel.ActiveWorkbook

# Update the view/selection on "gdp_calibrate" first so that it does not
# remain the active sheet/tab once we are done (the "config" sheet should
# stay the active tab, as in the original workbook).
$wsCalibrate = $wb.Worksheets.Item("gdp_calibrate")
[void]$wsCalibrate.Range("B3:B5").Select()

# Extend the "config" sheet with a new "year" column.
$wsConfig = $wb.Worksheets.Item("config")

$wsConfig.Range("E1").Value = "year"
$wsConfig.Range("E2").Value = 2020
$wsConfig.Range("E3").Value = 2030
$wsConfig.Range("E4").Value = 2040

[void]$wsConfig.Activate()
[void]$wsConfig.Range("E2:E4").Select()
